# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the refreshed scrape output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    6  = 12489
    7  = 12489
    11 = 452
    13 = 933
    14 = 13635
    15 = 13909
    23 = 426
    24 = 5010
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
